$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.247803211212158
$ws.Range("B1").Value = 2.665286302566528
$ws.Range("C1").Value = 8.439303398132324
$ws.Range("D1").Value = 2.10633659362793
$ws.Range("E1").Value = 1.141396999359131
